$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "21.01.2026, 1, 08:00-08:45, sala: 38"
$ws.Range("B5").Value = "21.01.2026, 8, 14:05-14:50, sala: 37"
$ws.Range("C5").Value = "Najwer Maciej"
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "3TH|JA2"
$ws.Range("F5").Value = "Informatyka"
$ws.Range("G5").Value = "informatyka, przeniesiona z lekcji 1 na lekcję 8"

$ws.Columns.Item(7).ColumnWidth = 41
